$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A4 with the new combined values
$ws.Range("A2").Value = "('Elemental', ['Token Creature — Elemental', 'Trample', '7/7'])"
$ws.Range("A3").Value = "('Elf Warrior', ['Token Creature — Elf Warrior', '1/1'])"
$ws.Range("A4").Value = "('Goblin', ['Token Creature — Goblin', '1/1'])"

# Remove the now-unused rows 5-11 so the sheet dimension shrinks to A1:A4
$ws.Range("A5:A11").EntireRow.Delete() | Out-Null
